# Automatic update of files.
# Update the "Förändrad" (changed) date column (C) for rows 2-19 from
# 45224 (2023-10-25) to 45233 (2023-11-03) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 19; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45224) {
        $cell.Value = 45233
    }
}
